# Todo.xlsx edit: "add create list and add item to list"
#
# Adds two new backlog rows (14 & 15) to the Client section, updates the
# sheet's default font from Arial to Calibri, and follows through on the
# knock-on formatting effects that font change has in Excel: the two
# bestFit columns get re-measured, the wrapped-text rows grow taller,
# and the selection lands on the next empty row below the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Default font: Arial -> Calibri ------------------------------------
# Applied to the existing used range only (so no empty cells get
# materialised); the new rows get the same treatment further down once
# they exist.
$ws.Range("A1:C10").Font.Name = "Calibri"

# --- New rows: "add create list and add item to list" ------------------
$ws.Range("A14").Value = "show head line of list"
$ws.Range("A15").Value = "view to show a list of lists in preview mode for select a list for edit."
$ws.Range("A14").Font.Name = "Calibri"
$ws.Range("A15").Font.Name = "Calibri"

# Row 15's text wraps to two lines at the current column width (same as
# row 7 already does), so its row height doubles, same as row 7 after
# the default row height grows from 14.25 (Arial) to 15 (Calibri).
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30

# --- Column widths: bestFit columns re-measured under the new font -----
$ws.Columns.Item(1).ColumnWidth = 38.570870535714285
$ws.Columns.Item(3).ColumnWidth = 107.71149553571429

# --- Selection moves to the next empty row below the new content -------
$ws.Range("A16").Select() | Out-Null
